# "Started sweet menu animation"
# Adds two new slides (Title and Content layout) right after the existing
# slide: "Game Focus" and "Game Mode", describing the game's input-focus
# routing and the various tube/game modes.

$p = $ppt.ActivePresentation

# --- Slide 2: "Game Focus" -------------------------------------------------
$s2 = $p.Slides.Add(2, 2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Game Focus"

$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange
$tr2.Text = "Controls Keyboard/Touch Events"
[void]$tr2.InsertAfter("`rFOCUS_MENU – routes input events to menu")
[void]$tr2.InsertAfter("`rFOCUS_TOWER – routes input events to tower")
$tr2.Paragraphs(2).IndentLevel = 2
$tr2.Paragraphs(3).IndentLevel = 2

# --- Slide 3: "Game Mode" ---------------------------------------------------
$s3 = $p.Slides.Add(3, 2)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Game Mode"

$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange
$tr3.Text = "Loading – Nothing Visible, game is loading everything"
[void]$tr3.InsertAfter("`rNone – Tube goes off screen, nothing visible.")
[void]$tr3.InsertAfter("`rClosed – Tube is closed, used for transitions. – Everything is visible.")
[void]$tr3.InsertAfter("`rEndless – standard endless")
[void]$tr3.InsertAfter("`rLine Clear – ")
[void]$tr3.InsertAfter("line clear game")
[void]$tr3.InsertAfter("`r")
[void]$tr3.InsertAfter("`r")
